$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-9 (A:T) with new TPM-derived values for the Gal-Galr3 LR pair table.
# New sending cluster "ECs" is introduced; row 4 (previously Resolving-Mac) now becomes FAPs-only,
# and the table now has 8 data rows (4 sending clusters x 2 target clusters).

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gal"
$ws.Cells.Item(2,3).Value = "Galr3"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.061714666666667
$ws.Cells.Item(2,8).Value = 3.185144
$ws.Cells.Item(2,9).Value = 0.09827129933343294
$ws.Cells.Item(2,10).Value = 0.09827129933343293
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.2113103333333334
$ws.Cells.Item(2,14).Value = 0.633931
$ws.Cells.Item(2,15).Value = 0.5841777494360321
$ws.Cells.Item(2,16).Value = 0.5841777494360321
$ws.Cells.Item(2,17).Value = 0.2243512801182223
$ws.Cells.Item(2,18).Value = 2.019161521064
$ws.Cells.Item(2,19).Value = 0.05740790647875949
$ws.Cells.Item(2,20).Value = 0.05740790647875949

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gal"
$ws.Cells.Item(3,3).Value = "Galr3"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.061714666666667
$ws.Cells.Item(3,8).Value = 3.185144
$ws.Cells.Item(3,9).Value = 0.09827129933343294
$ws.Cells.Item(3,10).Value = 0.09827129933343293
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.1504123333333333
$ws.Cells.Item(3,14).Value = 0.451237
$ws.Cells.Item(3,15).Value = 0.415822250563968
$ws.Cells.Item(3,16).Value = 0.415822250563968
$ws.Cells.Item(3,17).Value = 0.1596949803475556
$ws.Cells.Item(3,18).Value = 1.437254823128
$ws.Cells.Item(3,19).Value = 0.04086339285467346
$ws.Cells.Item(3,20).Value = 0.04086339285467345

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gal"
$ws.Cells.Item(4,3).Value = "Galr3"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.642758
$ws.Cells.Item(4,8).Value = 4.928274
$ws.Cells.Item(4,9).Value = 0.1520521174085614
$ws.Cells.Item(4,10).Value = 0.1520521174085614
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.2113103333333334
$ws.Cells.Item(4,14).Value = 0.633931
$ws.Cells.Item(4,15).Value = 0.5841777494360321
$ws.Cells.Item(4,16).Value = 0.5841777494360321
$ws.Cells.Item(4,17).Value = 0.347131740566
$ws.Cells.Item(4,18).Value = 3.124185665094
$ws.Cells.Item(4,19).Value = 0.08882546374471668
$ws.Cells.Item(4,20).Value = 0.08882546374471668

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gal"
$ws.Cells.Item(5,3).Value = "Galr3"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.642758
$ws.Cells.Item(5,8).Value = 4.928274
$ws.Cells.Item(5,9).Value = 0.1520521174085614
$ws.Cells.Item(5,10).Value = 0.1520521174085614
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.1504123333333333
$ws.Cells.Item(5,14).Value = 0.451237
$ws.Cells.Item(5,15).Value = 0.415822250563968
$ws.Cells.Item(5,16).Value = 0.415822250563968
$ws.Cells.Item(5,17).Value = 0.247091063882
$ws.Cells.Item(5,18).Value = 2.223819574938
$ws.Cells.Item(5,19).Value = 0.06322665366384468
$ws.Cells.Item(5,20).Value = 0.06322665366384468

# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Gal"
$ws.Cells.Item(6,3).Value = "Galr3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 7.792831666666667
$ws.Cells.Item(6,8).Value = 23.378495
$ws.Cells.Item(6,9).Value = 0.7212970842480482
$ws.Cells.Item(6,10).Value = 0.7212970842480481
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.2113103333333334
$ws.Cells.Item(6,14).Value = 0.633931
$ws.Cells.Item(6,15).Value = 0.5841777494360321
$ws.Cells.Item(6,16).Value = 0.5841777494360321
$ws.Cells.Item(6,17).Value = 1.646705857093889
$ws.Cells.Item(6,18).Value = 14.820352713845
$ws.Cells.Item(6,19).Value = 0.4213657073507968
$ws.Cells.Item(6,20).Value = 0.4213657073507968

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Gal"
$ws.Cells.Item(7,3).Value = "Galr3"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 7.792831666666667
$ws.Cells.Item(7,8).Value = 23.378495
$ws.Cells.Item(7,9).Value = 0.7212970842480482
$ws.Cells.Item(7,10).Value = 0.7212970842480481
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.1504123333333333
$ws.Cells.Item(7,14).Value = 0.451237
$ws.Cells.Item(7,15).Value = 0.415822250563968
$ws.Cells.Item(7,16).Value = 0.415822250563968
$ws.Cells.Item(7,17).Value = 1.172137994257222
$ws.Cells.Item(7,18).Value = 10.549241948315
$ws.Cells.Item(7,19).Value = 0.2999313768972514
$ws.Cells.Item(7,20).Value = 0.2999313768972514

# Row 8
$ws.Cells.Item(8,1).Value = "Resolving-Mac"
$ws.Cells.Item(8,2).Value = "Gal"
$ws.Cells.Item(8,3).Value = "Galr3"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.3066096666666667
$ws.Cells.Item(8,8).Value = 0.919829
$ws.Cells.Item(8,9).Value = 0.02837949900995756
$ws.Cells.Item(8,10).Value = 0.02837949900995756
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.2113103333333334
$ws.Cells.Item(8,14).Value = 0.633931
$ws.Cells.Item(8,15).Value = 0.5841777494360321
$ws.Cells.Item(8,16).Value = 0.5841777494360321
$ws.Cells.Item(8,17).Value = 0.06478979086655556
$ws.Cells.Item(8,18).Value = 0.5831081177990001
$ws.Cells.Item(8,19).Value = 0.01657867186175911
$ws.Cells.Item(8,20).Value = 0.01657867186175911

# Row 9
$ws.Cells.Item(9,1).Value = "Resolving-Mac"
$ws.Cells.Item(9,2).Value = "Gal"
$ws.Cells.Item(9,3).Value = "Galr3"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.3066096666666667
$ws.Cells.Item(9,8).Value = 0.919829
$ws.Cells.Item(9,9).Value = 0.02837949900995756
$ws.Cells.Item(9,10).Value = 0.02837949900995756
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.1504123333333333
$ws.Cells.Item(9,14).Value = 0.451237
$ws.Cells.Item(9,15).Value = 0.415822250563968
$ws.Cells.Item(9,16).Value = 0.415822250563968
$ws.Cells.Item(9,17).Value = 0.04611787538588889
$ws.Cells.Item(9,18).Value = 0.415060878473
$ws.Cells.Item(9,19).Value = 0.01180082714819846
$ws.Cells.Item(9,20).Value = 0.01180082714819846

Write-Output "Updated sheet1 rows 2-9 with new TPM data (dimension now A1:T9)."